$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows (15 & 20 agent experiments, 500 exp) to append after existing row 13
$data = @(
    @(15, 1, 62.742, 940.376, 0.93306206, 83036265.94400001, 0.06503414, 0.27340972, 17.97428276, 13.09818923066255, 196.3738271358175, 0.1939020838872756, 37308077.89811619, 0.01780541457638304, 0.09981963582438994, 9.776680470392529, 0.15),
    @(15, 1, 80.84, 1185.924, 0.73006058, 65727917.896, 0.05246292, 0.28368434, 23.7410607, 14.33154514952053, 205.9520106368732, 0.1236145603110953, 22546383.93325812, 0.005707393464087422, 0.06415970915709208, 9.673066310122771, 0.85),
    @(20, 1, 46.092, 920.924, 0.96314906, 123300097.872, 0.08834232, 0.2065659, 10.2250882, 10.87345391308576, 217.2577088346837, 0.2226476197854103, 57181264.09824292, 0.0293864971687559, 0.09457383351944199, 7.111083349348784, 0.15),
    @(20, 1, 56.512, 1098.198, 0.78682602, 100206206.928, 0.08038532, 0.22797276, 13.32970388, 9.887030025515797, 186.8924406302445, 0.1274289746325828, 30305669.43869139, 0.007630316241775239, 0.04991430324184279, 5.483287335673598, 0.85)
)

$startRow = 14
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    for ($c = 1; $c -le $data[$i].Count; $c++) {
        $ws.Cells.Item($row, $c).Value = $data[$i][$c - 1]
    }
}
